$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one record per row (rows 2..102). The edit inserts a brand
# new record at row 31 (pushing every following row down by one, so the
# former row 102 ends up at row 103) and fills the new row 31 with its own
# data (same Mercado/Region/etc. as its neighbour, but a new date/volume/
# price set).

# Insert a new row before row 31; existing row 31..102 shift down to 32..103.
$ws.Rows.Item(31).Insert(-4121, 0)

# Populate the newly inserted row 31 with the new record's values.
$ws.Cells.Item(31, 1).Value = 10
$ws.Cells.Item(31, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(31, 3).Value = "La Araucanía"
$ws.Cells.Item(31, 4).Value = Get-Date -Year 2023 -Month 4 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(31, 5).Value = 9
$ws.Cells.Item(31, 6).Value = 300000001
$ws.Cells.Item(31, 7).Value = "Rabanito"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 50
$ws.Cells.Item(31, 11).Value = 8000
$ws.Cells.Item(31, 12).Value = 8000
$ws.Cells.Item(31, 13).Value = 8000
$ws.Cells.Item(31, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(31, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(31, 16).Value = 667
$ws.Cells.Item(31, 17).Value = 12
$ws.Cells.Item(31, 18).Value = "Hortaliza"

# Match the date-cell number format used by the rest of column D.
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(32, 4).NumberFormat
